# "Mise a jour site" - add a "disponible" availability column (F) to the
# price list sheet: bold centered header in F1, and a centered "1" in
# F2:F72 for every existing product row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# Header cell
$header = $ws.Range("F1")
$header.Value = "disponible"
$header.Font.Bold = $true
$header.HorizontalAlignment = $xlCenter
$header.VerticalAlignment = $xlCenter

# Data rows: mark every product as available (1), centered.
$lastRow = 72
$data = $ws.Range("F2:F" + $lastRow)
$data.HorizontalAlignment = $xlCenter
$data.VerticalAlignment = $xlCenter

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = 1
}

# Match the saved selection/view state from the edit.
[void]$ws.Range("H6").Select()
